# edit.ps1 - applies the "Small adjustment for hololens section" edit.
#
# Summary of the three textual changes (per the target diff):
#   1. In the "HoloLens" intro paragraph, the sentence
#        "...Microsofts HoloLens verwendet. "
#      gets the phrase "auf Basis von Unity3D und C# " inserted between
#      "HoloLens " and "verwendet. ", ending up as three separate runs.
#   2. In the "Recommender" paragraph, the long run starting with
#      "-Informationen werden in einem kleinen Textfeld in einer Ecke..."
#      is split in two, with the (pre-existing) "_GoBack" bookmark
#      inserted exactly at the split point (after "...Textfeld in").
#   3. The "_GoBack" bookmark's old location (end of the next paragraph,
#      after "...gestartet werden muss.") is removed - which happens
#      automatically because Word bookmark names are unique, so re-adding
#      "_GoBack" at the new spot simply relocates it.
#
# Technique used to force Word to keep freshly-typed / freshly-split text
# in its own <w:r> run instead of silently re-merging it into a
# neighbouring run with identical formatting: drop a temporary bookmark
# exactly on the boundary we want to protect, then delete that bookmark
# via the Bookmarks collection. The bookmark tags disappear again, but
# the run split they forced survives the subsequent save/normalize step.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Microsofts HoloLens verwendet." -> insert
# "auf Basis von Unity3D und C# " between "HoloLens " and "verwendet. "
# ---------------------------------------------------------------------

$find1 = $d.Content
$find1.Find.Execute("HoloLens verwendet. ")
$s1 = $find1.Start
$e1 = $find1.End
$splitPoint1 = $s1 + 9   # right after "HoloLens " (9 chars incl. trailing space)

# Protect the boundary that follows this whole sentence (between
# "verwendet. " and the next run "Durch diese ") *before* we touch
# anything, so that run keeps its own identity/attributes.
$d.Bookmarks.Add("TempGuard1", $d.Range($e1, $e1))

# Insert the new phrase; at this point it merges into one run with
# whatever is on its left ("...HoloLens ").
$insertionPoint = $d.Range($splitPoint1, $splitPoint1)
$newPhrase = "auf Basis von Unity3D und C# "
$insertionPoint.InsertBefore($newPhrase)

# Now drop bookmarks on both edges of the freshly inserted phrase and
# delete them again - this peels the phrase off into its own run and
# separates it from "verwendet. " as well.
$boundaryBefore = $splitPoint1
$boundaryAfter = $splitPoint1 + $newPhrase.Length
$d.Bookmarks.Add("TempGuard2", $d.Range($boundaryBefore, $boundaryBefore))
$d.Bookmarks.Add("TempGuard3", $d.Range($boundaryAfter, $boundaryAfter))

$d.Bookmarks.Item("TempGuard2").Delete()
$d.Bookmarks.Item("TempGuard3").Delete()
$d.Bookmarks.Item("TempGuard1").Delete()

# ---------------------------------------------------------------------
# Change 2 + 3: split the "Recommender-Informationen..." run after
# "...Textfeld in" and drop the (moved) "_GoBack" bookmark there. Since
# bookmark names are unique, this automatically removes it from its old
# location at the end of the following paragraph.
# ---------------------------------------------------------------------

$find2 = $d.Content
$find2.Find.Execute("-Informationen werden in einem kleinen Textfeld in")
$splitPoint2 = $find2.End

$find3 = $d.Content
$find3.Find.Execute("server-seitig entschieden.")
$farEnd2 = $find3.End

# Protect the end of this long run (boundary with " Da die Exponat-...")
# before splitting, so that following run keeps its own identity.
$d.Bookmarks.Add("TempGuard4", $d.Range($farEnd2, $farEnd2))

# Placing the real "_GoBack" bookmark here both forces the run split and
# relocates the bookmark permanently (removing it from its old spot).
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint2, $splitPoint2))

$d.Bookmarks.Item("TempGuard4").Delete()
